# Populate the "Code" column (C) with the currency exchange rate codes.
# Values are written in the order C4, C3, C2 so that the resulting shared
# string table lists them as: 0,9852 / 117,6594 / 105,3540 (matching the
# order in which the new rows were originally authored).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "0,9852"
$ws.Range("C3").Value = "117,6594"
$ws.Range("C2").Value = "105,3540"
